# =============================================================================
# edit.ps1
#
# Updates the weekly NYPD 6th Precinct CompStat report:
#   - bumps the report "Volume/Number" label and the reporting week date
#     range shown in the header narrative text
#   - refreshes the weekly / 28-day / year-to-date / 2-year crime count
#     table (rows 15-33) with newly collected figures and their derived
#     percent-change columns
# =============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Header narrative text (rich-text shared strings).
#    All runs in each string share identical formatting, so replacing the
#    whole cell text preserves the visual appearance.
# -----------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"


# -----------------------------------------------------------------------
# 2) Cells that flip between the "not applicable" text marker (rendered
#    as "0" or "***.*") and a real numeric figure (or vice versa).
#    A plain .Value assignment would leave the previous cell style in
#    place, so the number-format is first transplanted from a donor cell
#    that already carries the desired style, and only then is the final
#    value (or text) applied.
# -----------------------------------------------------------------------
$ws.Cells.Item(15, 9).Copy() | Out-Null
$ws.Cells.Item(17, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> C17
$ws.Cells.Item(17, 3).Value = 2   # C17

$ws.Cells.Item(15, 9).Copy() | Out-Null
$ws.Cells.Item(20, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> D20
$ws.Cells.Item(20, 4).Value = 1   # D20

$ws.Cells.Item(15, 11).Copy() | Out-Null
$ws.Cells.Item(20, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> E20
$ws.Cells.Item(20, 5).Value = -100   # E20

$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(20, 6).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> F20
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(20, 6).PasteSpecial(-4163) | Out-Null   # xlPasteValues -> F20 (text "0")

$ws.Cells.Item(15, 9).Copy() | Out-Null
$ws.Cells.Item(33, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> D33
$ws.Cells.Item(33, 4).Value = 1   # D33

$ws.Cells.Item(15, 11).Copy() | Out-Null
$ws.Cells.Item(33, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> E33
$ws.Cells.Item(33, 5).Value = -100   # E33

$ws.Cells.Item(15, 9).Copy() | Out-Null
$ws.Cells.Item(33, 7).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> G33
$ws.Cells.Item(33, 7).Value = 1   # G33

$ws.Cells.Item(15, 11).Copy() | Out-Null
$ws.Cells.Item(33, 8).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> H33
$ws.Cells.Item(33, 8).Value = -100   # H33

$ws.Cells.Item(15, 9).Copy() | Out-Null
$ws.Cells.Item(33, 10).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> J33
$ws.Cells.Item(33, 10).Value = 1   # J33

$ws.Cells.Item(15, 11).Copy() | Out-Null
$ws.Cells.Item(33, 11).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> K33
$ws.Cells.Item(33, 11).Value = -100   # K33

$excel.CutCopyMode = $false

# -----------------------------------------------------------------------
# 3) Remaining cells keep their existing style; only the value changes.
# -----------------------------------------------------------------------
# Row 15
$ws.Cells.Item(15, 6).Value = 1   # F15

# Row 16
$ws.Cells.Item(16, 3).Value = 7   # C16
$ws.Cells.Item(16, 4).Value = 5   # D16
$ws.Cells.Item(16, 5).Value = 40   # E16
$ws.Cells.Item(16, 6).Value = 16   # F16
$ws.Cells.Item(16, 7).Value = 12   # G16
$ws.Cells.Item(16, 8).Value = 33.333333333333   # H16
$ws.Cells.Item(16, 9).Value = 40   # I16
$ws.Cells.Item(16, 10).Value = 61   # J16
$ws.Cells.Item(16, 11).Value = -34.426229508196   # K16
$ws.Cells.Item(16, 12).Value = -38.461538461538   # L16
$ws.Cells.Item(16, 13).Value = -18.367346938775   # M16
$ws.Cells.Item(16, 14).Value = -87.577639751552   # N16

# Row 17
$ws.Cells.Item(17, 4).Value = 4   # D17
$ws.Cells.Item(17, 5).Value = -50   # E17
$ws.Cells.Item(17, 6).Value = 9   # F17
$ws.Cells.Item(17, 8).Value = -43.75   # H17
$ws.Cells.Item(17, 9).Value = 39   # I17
$ws.Cells.Item(17, 10).Value = 44   # J17
$ws.Cells.Item(17, 11).Value = -11.363636363636   # K17
$ws.Cells.Item(17, 12).Value = -37.096774193548   # L17
$ws.Cells.Item(17, 13).Value = 8.333333333333   # M17
$ws.Cells.Item(17, 14).Value = -56.666666666666   # N17

# Row 18
$ws.Cells.Item(18, 3).Value = 4   # C18
$ws.Cells.Item(18, 4).Value = 5   # D18
$ws.Cells.Item(18, 5).Value = -20   # E18
$ws.Cells.Item(18, 6).Value = 14   # F18
$ws.Cells.Item(18, 7).Value = 33   # G18
$ws.Cells.Item(18, 8).Value = -57.575757575757   # H18
$ws.Cells.Item(18, 9).Value = 59   # I18
$ws.Cells.Item(18, 10).Value = 98   # J18
$ws.Cells.Item(18, 11).Value = -39.795918367346   # K18
$ws.Cells.Item(18, 12).Value = -51.639344262295   # L18
$ws.Cells.Item(18, 13).Value = -25.316455696202   # M18
$ws.Cells.Item(18, 14).Value = -78.853046594982   # N18

# Row 19
$ws.Cells.Item(19, 3).Value = 15   # C19
$ws.Cells.Item(19, 4).Value = 21   # D19
$ws.Cells.Item(19, 5).Value = -28.571428571428   # E19
$ws.Cells.Item(19, 6).Value = 72   # F19
$ws.Cells.Item(19, 7).Value = 84   # G19
$ws.Cells.Item(19, 8).Value = -14.285714285714   # H19
$ws.Cells.Item(19, 9).Value = 336   # I19
$ws.Cells.Item(19, 10).Value = 376   # J19
$ws.Cells.Item(19, 11).Value = -10.638297872340   # K19
$ws.Cells.Item(19, 12).Value = -26.477024070021   # L19
$ws.Cells.Item(19, 13).Value = -8.695652173913   # M19
$ws.Cells.Item(19, 14).Value = -60.701754385964   # N19

# Row 20
$ws.Cells.Item(20, 7).Value = 4   # G20
$ws.Cells.Item(20, 8).Value = -100   # H20
$ws.Cells.Item(20, 10).Value = 17   # J20
$ws.Cells.Item(20, 11).Value = -70.588235294117   # K20
$ws.Cells.Item(20, 13).Value = -61.538461538461   # M20
$ws.Cells.Item(20, 14).Value = -98.293515358361   # N20

# Row 21
$ws.Cells.Item(21, 3).Value = 28   # C21
$ws.Cells.Item(21, 4).Value = 36   # D21
$ws.Cells.Item(21, 5).Value = -22.222222222222   # E21
$ws.Cells.Item(21, 6).Value = 112   # F21
$ws.Cells.Item(21, 7).Value = 149   # G21
$ws.Cells.Item(21, 8).Value = -24.832214765100   # H21
$ws.Cells.Item(21, 9).Value = 484   # I21
$ws.Cells.Item(21, 10).Value = 597   # J21
$ws.Cells.Item(21, 11).Value = -18.92797319933   # K21
$ws.Cells.Item(21, 12).Value = -32.871012482663   # L21
$ws.Cells.Item(21, 13).Value = -11.83970856102   # M21
$ws.Cells.Item(21, 14).Value = -73.752711496746   # N21

# Row 22
$ws.Cells.Item(22, 6).Value = 2   # F22
$ws.Cells.Item(22, 9).Value = 21   # I22
$ws.Cells.Item(22, 11).Value = 31.25   # K22
$ws.Cells.Item(22, 12).Value = 16.666666666666   # L22
$ws.Cells.Item(22, 13).Value = -16   # M22

# Row 24
$ws.Cells.Item(24, 3).Value = 35   # C24
$ws.Cells.Item(24, 4).Value = 27   # D24
$ws.Cells.Item(24, 5).Value = 29.629629629629   # E24
$ws.Cells.Item(24, 6).Value = 118   # F24
$ws.Cells.Item(24, 7).Value = 129   # G24
$ws.Cells.Item(24, 8).Value = -8.527131782945   # H24
$ws.Cells.Item(24, 9).Value = 549   # I24
$ws.Cells.Item(24, 10).Value = 639   # J24
$ws.Cells.Item(24, 11).Value = -14.084507042253   # K24
$ws.Cells.Item(24, 12).Value = -17.814371257485   # L24
$ws.Cells.Item(24, 13).Value = -4.188481675392   # M24

# Row 25
$ws.Cells.Item(25, 3).Value = 33   # C25
$ws.Cells.Item(25, 4).Value = 26   # D25
$ws.Cells.Item(25, 5).Value = 26.923076923076   # E25
$ws.Cells.Item(25, 6).Value = 81   # F25
$ws.Cells.Item(25, 7).Value = 99   # G25
$ws.Cells.Item(25, 8).Value = -18.181818181818   # H25
$ws.Cells.Item(25, 9).Value = 397   # I25
$ws.Cells.Item(25, 10).Value = 525   # J25
$ws.Cells.Item(25, 11).Value = -24.380952380952   # K25
$ws.Cells.Item(25, 12).Value = -17.291666666666   # L25

# Row 26
$ws.Cells.Item(26, 3).Value = 8   # C26
$ws.Cells.Item(26, 4).Value = 5   # D26
$ws.Cells.Item(26, 5).Value = 60   # E26
$ws.Cells.Item(26, 6).Value = 34   # F26
$ws.Cells.Item(26, 7).Value = 23   # G26
$ws.Cells.Item(26, 8).Value = 47.826086956521   # H26
$ws.Cells.Item(26, 9).Value = 133   # I26
$ws.Cells.Item(26, 10).Value = 133   # J26
$ws.Cells.Item(26, 11).Value = 0   # K26
$ws.Cells.Item(26, 12).Value = -16.875   # L26
$ws.Cells.Item(26, 13).Value = 47.777777777777   # M26

# Row 27
$ws.Cells.Item(27, 6).Value = 1   # F27

# Row 28
$ws.Cells.Item(28, 3).Value = 2   # C28
$ws.Cells.Item(28, 5).Value = 100   # E28
$ws.Cells.Item(28, 7).Value = 3   # G28
$ws.Cells.Item(28, 8).Value = 133.333333333333   # H28
$ws.Cells.Item(28, 9).Value = 31   # I28
$ws.Cells.Item(28, 10).Value = 23   # J28
$ws.Cells.Item(28, 11).Value = 34.782608695652   # K28
$ws.Cells.Item(28, 12).Value = 14.814814814814   # L28

# Row 31
$ws.Cells.Item(31, 6).Value = 1   # F31
$ws.Cells.Item(31, 7).Value = 5   # G31
$ws.Cells.Item(31, 8).Value = -80   # H31
$ws.Cells.Item(31, 10).Value = 10   # J31
$ws.Cells.Item(31, 11).Value = -40   # K31
$ws.Cells.Item(31, 12).Value = 200   # L31

